# Auto-generated Excel COM-interop script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on column D so numeric-looking strings (e.g. "1.00", "42.00",
# "0.0000175") are not silently coerced to Excel numbers and lose formatting.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.641.57"
$ws.Range("D3").Value = "3.451.30"
$ws.Range("E3").Value = "  +2.78%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "578.59"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").Value = "147.78"
$ws.Range("E6").Value = "  +7.90%  "
$ws.Range("D7").Value = "3.450.97"
$ws.Range("E7").Value = "  +2.81%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").Value = "7.69"
$ws.Range("E10").Value = "  +2.79%  "
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").Value = "0.387"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").Value = "4.039.81"
$ws.Range("E13").Value = "  +2.76%  "
$ws.Range("D14").Value = "27.90"
$ws.Range("E14").Value = "  +7.27%  "
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "0.0000175"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").Value = "3.449.77"
$ws.Range("E17").Value = "  +2.69%  "
$ws.Range("D18").Value = "61.680.34"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "6.27"
$ws.Range("E19").Value = "  +7.20%  "
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").Value = "9.43"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").Value = "384.03"
$ws.Range("E22").Value = "  +1.56%  "
$ws.Range("E23").Value = "  +2.56%  "
$ws.Range("D24").Value = "3.589.07"
$ws.Range("E24").Value = "  +2.56%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").Value = "72.23"
$ws.Range("E27").Value = "  +1.38%  "
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("D29").Value = "0.177"
$ws.Range("E29").Value = "  +8.20%  "
$ws.Range("D30").Value = "7.77"
$ws.Range("E30").Value = "  +3.98%  "
$ws.Range("D31").Value = "1.57"
$ws.Range("E31").Value = "  -13.34%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D34").Value = "2.17"
$ws.Range("E34").Value = "  +1.14%  "
$ws.Range("D36").Value = "24.05"
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("D37").Value = "5.23"
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("D38").Value = "7.02"
$ws.Range("E38").Value = "  +3.26%  "
$ws.Range("E39").Value = "  +2.26%  "
$ws.Range("D40").Value = "166.21"
$ws.Range("D41").Value = "0.0784"
$ws.Range("E41").Value = "  +2.71%  "
$ws.Range("D42").Value = "25.88"
$ws.Range("E42").Value = "  +8.60%  "
$ws.Range("E43").Value = "  +2.91%  "
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("D47").Value = "42.00"
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "1.18"
$ws.Range("E48").Value = "  -2.28%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.617.38"
$ws.Range("E49").Value = "  +10.32%  "
$ws.Range("D50").Value = "23.64"
$ws.Range("E50").Value = "  +3.49%  "
$ws.Range("E51").Value = "  +0.48%  "
